$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Range("B4").Value = "In Translation"
$ws1.Range("C4").Value = "In Translation"

$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Range("B4").Value = "In Translation"

$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Range("B4").Value = "In Translation"
